$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update parameter values (B2:B13)
$ws.Range("B2").Value = 0.1089559434280109
$ws.Range("B3").Value = -0.0179798733778174
$ws.Range("B4").Value = -0.1757637317637759
$ws.Range("B5").Value = 0.0710079908073632
$ws.Range("B6").Value = -0.380935120180237
$ws.Range("B7").Value = 0.1458385850570895
$ws.Range("B8").Value = 0.194368300192797
$ws.Range("B9").Value = -0.001626525494170596
$ws.Range("B10").Value = -0.744767904920193
$ws.Range("B11").Value = -0.6697511584247826
$ws.Range("B12").Value = -0.6476415368041971
$ws.Range("B13").Value = 0.3112709005875231

# Remove row 14 (shot_during_regular_play) entirely
$ws.Rows.Item(14).Delete()
